# BISAGRAS FICHAS carpintero.xlsx - "fix bug exeded requeste in google drive"
#
# The sheet header date (A1) is bumped by one day, and the three hinge
# prices (D29:D31) are updated to their new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Date stamp in the header (serial 45310 -> 45311, i.e. 2024-01-19 -> 2024-01-20)
$ws.Range("A1").Value = 45311

# Updated hinge prices
$ws.Range("D29").Value = 201.233
$ws.Range("D30").Value = 250
$ws.Range("D31").Value = 229.232
